$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 6 new data rows before the old last row (old row 28) ---
# This pushes the old "closing-style" row 28 down to row 34, and the
# inserted rows 28-33 inherit formatting from the row above.
$ws.Rows.Item(28).Resize(6).Insert()

# Re-apply the "normal" data-row formatting (style of row 27) onto the
# newly inserted rows 28-33 so every row from 16-33 shares the same
# look, matching row 34's pre-existing "closing" border style.
$ws.Range("B27:J27").Copy($ws.Range("B28:J33"))

# --- 2. Update the summary / header fields ---
$ws.Range("E11").Value = 692064
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 19

# --- 3. Rewrite all 19 detail rows (16-34) with the new data set ---
$ws.Cells.Item(16,2).Value = "CC"
$ws.Cells.Item(16,3).Value = "1052734505"
$ws.Cells.Item(16,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(16,5).Value = "2112"
$ws.Cells.Item(16,6).Value = 35112
$ws.Cells.Item(16,7).Value = 908526
$ws.Cells.Item(17,2).Value = "CC"
$ws.Cells.Item(17,3).Value = "1052734505"
$ws.Cells.Item(17,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(17,5).Value = "2111"
$ws.Cells.Item(17,6).Value = 35112
$ws.Cells.Item(17,7).Value = 908526
$ws.Cells.Item(18,2).Value = "CC"
$ws.Cells.Item(18,3).Value = "1052734505"
$ws.Cells.Item(18,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(18,5).Value = "2110"
$ws.Cells.Item(18,6).Value = 35112
$ws.Cells.Item(18,7).Value = 908526
$ws.Cells.Item(19,2).Value = "CC"
$ws.Cells.Item(19,3).Value = "1052734505"
$ws.Cells.Item(19,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(19,5).Value = "2109"
$ws.Cells.Item(19,6).Value = 35112
$ws.Cells.Item(19,7).Value = 908526
$ws.Cells.Item(20,2).Value = "CC"
$ws.Cells.Item(20,3).Value = "1052734505"
$ws.Cells.Item(20,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(20,5).Value = "2108"
$ws.Cells.Item(20,6).Value = 35112
$ws.Cells.Item(20,7).Value = 908526
$ws.Cells.Item(21,2).Value = "CC"
$ws.Cells.Item(21,3).Value = "1052734505"
$ws.Cells.Item(21,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(21,5).Value = "2107"
$ws.Cells.Item(21,6).Value = 35112
$ws.Cells.Item(21,7).Value = 908526
$ws.Cells.Item(22,2).Value = "CC"
$ws.Cells.Item(22,3).Value = "1052734505"
$ws.Cells.Item(22,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(22,5).Value = "2106"
$ws.Cells.Item(22,6).Value = 35112
$ws.Cells.Item(22,7).Value = 908526
$ws.Cells.Item(23,2).Value = "CC"
$ws.Cells.Item(23,3).Value = "1052734505"
$ws.Cells.Item(23,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(23,5).Value = "2105"
$ws.Cells.Item(23,6).Value = 35112
$ws.Cells.Item(23,7).Value = 908526
$ws.Cells.Item(24,2).Value = "CC"
$ws.Cells.Item(24,3).Value = "1052734505"
$ws.Cells.Item(24,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(24,5).Value = "2104"
$ws.Cells.Item(24,6).Value = 35112
$ws.Cells.Item(24,7).Value = 908526
$ws.Cells.Item(25,2).Value = "CC"
$ws.Cells.Item(25,3).Value = "1052734505"
$ws.Cells.Item(25,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(25,5).Value = "2103"
$ws.Cells.Item(25,6).Value = 35112
$ws.Cells.Item(25,7).Value = 908526
$ws.Cells.Item(26,2).Value = "CC"
$ws.Cells.Item(26,3).Value = "1052734505"
$ws.Cells.Item(26,4).Value = "RAFAEL DIONICIO GARCIA SIERRA"
$ws.Cells.Item(26,5).Value = "2102"
$ws.Cells.Item(26,6).Value = 35112
$ws.Cells.Item(26,7).Value = 908526
$ws.Cells.Item(27,2).Value = "CC"
$ws.Cells.Item(27,3).Value = "1048604929"
$ws.Cells.Item(27,4).Value = "CRISTIAN DAVID LEYVA GUARDO"
$ws.Cells.Item(27,5).Value = "2101"
$ws.Cells.Item(27,6).Value = 32707
$ws.Cells.Item(27,7).Value = 908526
$ws.Cells.Item(28,2).Value = "CC"
$ws.Cells.Item(28,3).Value = "20191039"
$ws.Cells.Item(28,4).Value = "EVER ENRIQUE LEIVA REBOLLEDO"
$ws.Cells.Item(28,5).Value = "1902"
$ws.Cells.Item(28,6).Value = 33125
$ws.Cells.Item(28,7).Value = 828116
$ws.Cells.Item(29,2).Value = "CC"
$ws.Cells.Item(29,3).Value = "9201403"
$ws.Cells.Item(29,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(29,5).Value = "2507"
$ws.Cells.Item(29,6).Value = 40000
$ws.Cells.Item(29,7).Value = 1000000
$ws.Cells.Item(30,2).Value = "CC"
$ws.Cells.Item(30,3).Value = "9201403"
$ws.Cells.Item(30,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(30,5).Value = "2506"
$ws.Cells.Item(30,6).Value = 40000
$ws.Cells.Item(30,7).Value = 1000000
$ws.Cells.Item(31,2).Value = "CC"
$ws.Cells.Item(31,3).Value = "9201403"
$ws.Cells.Item(31,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(31,5).Value = "2505"
$ws.Cells.Item(31,6).Value = 40000
$ws.Cells.Item(31,7).Value = 1000000
$ws.Cells.Item(32,2).Value = "CC"
$ws.Cells.Item(32,3).Value = "9201403"
$ws.Cells.Item(32,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(32,5).Value = "2504"
$ws.Cells.Item(32,6).Value = 40000
$ws.Cells.Item(32,7).Value = 1000000
$ws.Cells.Item(33,2).Value = "CC"
$ws.Cells.Item(33,3).Value = "9201403"
$ws.Cells.Item(33,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(33,5).Value = "2503"
$ws.Cells.Item(33,6).Value = 40000
$ws.Cells.Item(33,7).Value = 1000000
$ws.Cells.Item(34,2).Value = "CC"
$ws.Cells.Item(34,3).Value = "9201403"
$ws.Cells.Item(34,4).Value = "MEDARDO PEREZ CABARCAS"
$ws.Cells.Item(34,5).Value = "2502"
$ws.Cells.Item(34,6).Value = 40000
$ws.Cells.Item(34,7).Value = 1000000

# --- 4. Clear any leftover formula/number formatting artifacts is not
# needed; header labels (row 15), footer labels (now rows 39-40) and
# merged cell ranges shift automatically with the row insert above.
